$d = $word.ActiveDocument
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $pf = $d.Paragraphs($i).Range.ParagraphFormat
    if ($pf.LeftIndent -gt 0 -or $pf.RightIndent -gt 0) {
        $pf.CharacterUnitFirstLineIndent = 0
    }
}
